$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 632.5
$ws.Range("I11").Value = 632.5
$ws.Range("K11").Value = 632.5
$ws.Range("M11").Value = -492.5

$ws.Range("H86").Value = 173615090
$ws.Range("I86").Value = 333335000
$ws.Range("K86").Value = 333335000
$ws.Range("M86").Value = -333333877

$ws.Range("H88").Value = 14529986
$ws.Range("I88").Value = 47622850
$ws.Range("J88").Value = 51858.688
$ws.Range("K88").Value = 47622850
$ws.Range("L88").Value = 51858.688
$ws.Range("M88").Value = -47622444
$ws.Range("N88").Value = -52670.688

$ws.Range("H89").Value = 173615090
$ws.Range("I89").Value = 333335000
$ws.Range("K89").Value = 1666675000
$ws.Range("M89").Value = -1666669384

$ws.Range("H91").Value = 14529986
$ws.Range("I91").Value = 47622850
$ws.Range("J91").Value = 51858.688
$ws.Range("K91").Value = 47622850
$ws.Range("L91").Value = 51858.688
$ws.Range("M91").Value = -47621446
$ws.Range("N91").Value = -54666.688

$ws.Range("H106").Value = 4998.6
$ws.Range("I106").Value = 4998.6
$ws.Range("K106").Value = 4998.6
$ws.Range("M106").Value = -4367.6

$ws.Range("H132").Value = 4562.8184
$ws.Range("I132").Value = 4025
$ws.Range("J132").Value = 5997
$ws.Range("K132").Value = 12075
$ws.Range("L132").Value = 17991
$ws.Range("M132").Value = -9545
$ws.Range("N132").Value = -23051

$ws.Range("H135").Value = 370968.75
$ws.Range("I135").Value = 435276.56
$ws.Range("J135").Value = 1198.75
$ws.Range("K135").Value = 3917489.04
$ws.Range("L135").Value = 10788.75
$ws.Range("M135").Value = -3914954.04
$ws.Range("N135").Value = -15858.75

$ws.Range("H138").Value = 7678.2
$ws.Range("I138").Value = 1740.5555
$ws.Range("J138").Value = 12536.272
$ws.Range("K138").Value = 5221.666499999999
$ws.Range("L138").Value = 37608.81600000001
$ws.Range("M138").Value = -81.66649999999936
$ws.Range("N138").Value = -47888.81600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3053.6897
$ws.Range("I2").Value = 1898.0555
$ws.Range("K2").Value = 1898.0555
$ws.Range("M2").Value = -1785.0555

$ws.Range("H32").Value = 2198066.8
$ws.Range("I32").Value = 2276950.8
$ws.Range("K32").Value = 2276950.8
$ws.Range("M32").Value = -2276663.8

$ws.Range("H61").Value = 7083.2563
$ws.Range("I61").Value = 3323
$ws.Range("K61").Value = 3323
$ws.Range("M61").Value = -3111

$ws.Range("H74").Value = 49709
$ws.Range("I74").Value = 65747
$ws.Range("K74").Value = 65747
$ws.Range("M74").Value = -64873

$ws.Range("H77").Value = 49709
$ws.Range("I77").Value = 65747
$ws.Range("K77").Value = 328735
$ws.Range("M77").Value = -324367

$ws.Range("H97").Value = 4395052
$ws.Range("I97").Value = 596.6
$ws.Range("J97").Value = 9277780
$ws.Range("K97").Value = 596.6
$ws.Range("L97").Value = 9277780
$ws.Range("M97").Value = -100.6
$ws.Range("N97").Value = -9278772

$ws.Range("H102").Value = 1198.7059
$ws.Range("I102").Value = 1255.5714
$ws.Range("J102").Value = 933.3333
$ws.Range("K102").Value = 1255.5714
$ws.Range("L102").Value = 933.3333
$ws.Range("M102").Value = 366.4286
$ws.Range("N102").Value = -4177.3333

$ws.Range("H116").Value = 3053.6897
$ws.Range("I116").Value = 1898.0555
$ws.Range("K116").Value = 1898.0555
$ws.Range("M116").Value = 395.9445000000001

$ws.Range("H136").Value = 7083.2563
$ws.Range("I136").Value = 3323
$ws.Range("K136").Value = 9969
$ws.Range("M136").Value = -7419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3053.6897
$ws.Range("I3").Value = 1898.0555
$ws.Range("K3").Value = 1898.0555
$ws.Range("M3").Value = -1784.0555

$ws.Range("H86").Value = 36767824
$ws.Range("I86").Value = 16669369
$ws.Range("K86").Value = 16669369
$ws.Range("M86").Value = -16668246

$ws.Range("H89").Value = 36767824
$ws.Range("I89").Value = 16669369
$ws.Range("K89").Value = 83346845
$ws.Range("M89").Value = -83341229

$ws.Range("H94").Value = 1536.7142
$ws.Range("I94").Value = 712.5454999999999
$ws.Range("J94").Value = 4558.6665
$ws.Range("K94").Value = 712.5454999999999
$ws.Range("L94").Value = 4558.6665
$ws.Range("M94").Value = -261.5454999999999
$ws.Range("N94").Value = -5460.6665

$ws.Range("H99").Value = 3032825
$ws.Range("I99").Value = 2136.6
$ws.Range("K99").Value = 2136.6
$ws.Range("M99").Value = -638.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10550.114
$ws.Range("I31").Value = 5002.5
$ws.Range("J31").Value = 14248.523
$ws.Range("K31").Value = 5002.5
$ws.Range("L31").Value = 14248.523
$ws.Range("M31").Value = -4707.5
$ws.Range("N31").Value = -14838.523

$ws.Range("H34").Value = 10550.114
$ws.Range("I34").Value = 5002.5
$ws.Range("J34").Value = 14248.523
$ws.Range("K34").Value = 5002.5
$ws.Range("L34").Value = 14248.523
$ws.Range("M34").Value = -4800.5
$ws.Range("N34").Value = -14652.523

$ws.Range("H132").Value = 7850.913
$ws.Range("I132").Value = 2668
$ws.Range("J132").Value = 10118.4375
$ws.Range("K132").Value = 8004
$ws.Range("L132").Value = 30355.3125
$ws.Range("M132").Value = -5474
$ws.Range("N132").Value = -35415.3125

$ws.Range("H134").Value = 9583.951999999999
$ws.Range("I134").Value = 2064.8
$ws.Range("K134").Value = 6194.400000000001
$ws.Range("M134").Value = -3659.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4960.375
$ws.Range("I56").Value = 4960.375
$ws.Range("K56").Value = 4960.375
$ws.Range("M56").Value = -4430.375

$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45856

$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47964

$ws.Range("H103").Value = 501.44446
$ws.Range("I103").Value = 348.5
$ws.Range("K103").Value = 1045.5
$ws.Range("M103").Value = -166.5

$ws.Range("H131").Value = 1798
$ws.Range("I131").Value = 747.375
$ws.Range("J131").Value = 3198.8333
$ws.Range("K131").Value = 2242.125
$ws.Range("L131").Value = 9596.499899999999
$ws.Range("M131").Value = 2797.875
$ws.Range("N131").Value = -19676.4999

$ws.Range("H141").Value = 9567.909
$ws.Range("I141").Value = 2541.3333
$ws.Range("J141").Value = 17999.8
$ws.Range("K141").Value = 7623.999899999999
$ws.Range("L141").Value = 53999.39999999999
$ws.Range("M141").Value = -2443.999899999999
$ws.Range("N141").Value = -64359.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4166.1816
$ws.Range("I80").Value = 3023.5
$ws.Range("K80").Value = 3023.5
$ws.Range("M80").Value = -2025.5

$ws.Range("H83").Value = 4166.1816
$ws.Range("I83").Value = 3023.5
$ws.Range("K83").Value = 15117.5
$ws.Range("M83").Value = -10125.5

$ws.Range("H102").Value = 2154.9678
$ws.Range("I102").Value = 1860.8334
$ws.Range("J102").Value = 3163.4285
$ws.Range("K102").Value = 1860.8334
$ws.Range("L102").Value = 3163.4285
$ws.Range("M102").Value = -238.8334
$ws.Range("N102").Value = -6407.4285

$ws.Range("H113").Value = 7778.5713
$ws.Range("I113").Value = 3928.5715
$ws.Range("K113").Value = 3928.5715
$ws.Range("M113").Value = -1758.5715

$ws.Range("H132").Value = 4441.9375
$ws.Range("I132").Value = 1760.9474
$ws.Range("J132").Value = 8360.308000000001
$ws.Range("K132").Value = 5282.8422
$ws.Range("L132").Value = 25080.924
$ws.Range("M132").Value = -2752.8422
$ws.Range("N132").Value = -30140.924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6989.4736
$ws.Range("I40").Value = 5374.5
$ws.Range("K40").Value = 5374.5
$ws.Range("M40").Value = -5238.5

$ws.Range("H82").Value = 126762920
$ws.Range("I82").Value = 202818300
$ws.Range("K82").Value = 202818300
$ws.Range("M82").Value = -202817939

$ws.Range("H85").Value = 126762920
$ws.Range("I85").Value = 202818300
$ws.Range("K85").Value = 202818300
$ws.Range("M85").Value = -202817052

$ws.Range("H122").Value = 6153.6665
$ws.Range("I122").Value = 3133.3333
$ws.Range("J122").Value = 8167.222
$ws.Range("K122").Value = 9399.999899999999
$ws.Range("L122").Value = 24501.666
$ws.Range("M122").Value = -6949.999899999999
$ws.Range("N122").Value = -29401.666

$ws.Range("H132").Value = 11911950
$ws.Range("I132").Value = 29414842
$ws.Range("J132").Value = 9983.76
$ws.Range("K132").Value = 88244526
$ws.Range("L132").Value = 29951.28
$ws.Range("M132").Value = -88241996
$ws.Range("N132").Value = -35011.28

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2541250
$ws.Range("I5").Value = 72500
$ws.Range("J5").Value = 5010000
$ws.Range("K5").Value = 72500
$ws.Range("L5").Value = 5010000
$ws.Range("M5").Value = -72388
$ws.Range("N5").Value = -5010224

$ws.Range("H15").Value = 24996.5
$ws.Range("I15").Value = 24996.5
$ws.Range("K15").Value = 24996.5
$ws.Range("M15").Value = -24708.5

$ws.Range("H54").Value = 14727.182
$ws.Range("J54").Value = 13499.5
$ws.Range("L54").Value = 13499.5
$ws.Range("N54").Value = -14539.5

$ws.Range("H81").Value = 9134864
$ws.Range("I81").Value = 1113177.4
$ws.Range("J81").Value = 14291663
$ws.Range("K81").Value = 2226354.8
$ws.Range("L81").Value = 28583326
$ws.Range("M81").Value = -2225293.8
$ws.Range("N81").Value = -28585448

$ws.Range("H84").Value = 9134864
$ws.Range("I84").Value = 1113177.4
$ws.Range("J84").Value = 14291663
$ws.Range("K84").Value = 11131774
$ws.Range("L84").Value = 142916630
$ws.Range("M84").Value = -11126470
$ws.Range("N84").Value = -142927238
